# BOM.xlsx update: uprev to 1.1 to support current parts availability;
# added solderpaste to USB connector through-pins.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Part renumbering: C15 -> C14 (Ref Des list for the 0.1uF cap row)
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "C1,C14,C25"

# ---------------------------------------------------------------------
# 2. Add solder paste cost to the USB connector through-pins row (row 17)
# ---------------------------------------------------------------------
$ws.Range("H17").Value = 0.39

# ---------------------------------------------------------------------
# 3. New BOM line: U1 FT230XS (USB to serial chip)
#    Values are written in this particular column order so that the
#    workbook's shared-string table is built up in the same order as
#    the authoritative edit (C before B, F before D).
# ---------------------------------------------------------------------
$ws.Range("A18").Value = "U1"
$ws.Range("C18").Value = "FT230XS"
$ws.Range("B18").Value = "USB to serial"
$ws.Range("F18").Value = "768-1135-1-ND"
$ws.Range("D18").Value = "16SSOP"
$ws.Range("E18").Value = "Digikey"
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 2.04

# ---------------------------------------------------------------------
# 4. Extend the "Cost badge" (column I) formula down through the new rows
# ---------------------------------------------------------------------
$ws.Range("I17:I18").Formula = "=G17*H17"

# ---------------------------------------------------------------------
# 5. Restore the author's last selection/cursor position
# ---------------------------------------------------------------------
$ws.Range("B16").Select()
